# Rewrite the 'Results' sheet data (rows 2-37) to match the target state, then
# drop the now-superfluous trailing row (old row 38).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Nvidia''s profits soar as AI boom shows no sign of slowing down ...'
$ws.Range("B2").Value = 45434
$ws.Range("C2").Value = '“Companies and countries are partnering with Nvidia to shift the trillion-dollar traditional data centres to accelerated computing and build'
$ws.Range("D2").Value = './Output/Images/0.jpg'
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = $false

$ws.Range("A3").Value = 'Breaking News, World News and Video from Al Jazeera'
$ws.Range("B3").Value = 45434
$ws.Range("C3").Value = 'News, analysis from the Middle East & worldwide, multimedia & interactives, opinions, documentaries, podcasts, long reads and broadcast'
$ws.Range("D3").Value = './Output/Images/1.jpg'
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = $false

$ws.Range("A4").Value = 'Israel''s war on Gaza live news: Attacks on besieged enclave kill 62 ...'
$ws.Range("B4").Value = 45434
$ws.Range("C4").Value = 'Israel, a major recipient of US military assistance for decades, is still due to receive billions of dollars of US aid and weaponry. “The'
$ws.Range("D4").Value = './Output/Images/2.jpg'
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = $false

$ws.Range("A5").Value = 'Economy | Today''s latest from Al Jazeera'
$ws.Range("B5").Value = 45433
$ws.Range("C5").Value = 'Russian court seizes two European banks'' assets amid Western sanctions. Freezing hundreds of billions of dollars in lenders'' assets was part of dispute over gas'
$ws.Range("D5").Value = './Output/Images/3.jpg'
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = $false

$ws.Range("A6").Value = 'Israel''s war on Gaza live news: Deadly combat rages as Rafah ...'
$ws.Range("B6").Value = 45434
$ws.Range("C6").Value = 'Israel, a major recipient of US military assistance for decades, is still due to receive billions of dollars of US aid and weaponry. “The'
$ws.Range("D6").Value = './Output/Images/4.jpg'
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = $false

$ws.Range("A7").Value = 'Gaza war: What does victory look like for the US and Israel? | Israel ...'
$ws.Range("B7").Value = 45433
$ws.Range("C7").Value = 'Israel has said it is seeking an “absolute victory” over Hamas, as it continues to receive billions of dollars in unconditional military aid'
$ws.Range("D7").Value = './Output/Images/5.jpg'
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = $false

$ws.Range("A8").Value = 'Israel''s war on Gaza live: UNRWA suspends food distribution in Rafah'
$ws.Range("B8").Value = 45432
$ws.Range("C8").Value = 'dollars in US military assistance that remains in the pipeline for Prime Minister Benjamin Netanyahu''s government. But Biden has also faced'
$ws.Range("D8").Value = './Output/Images/6.jpg'
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = $false

$ws.Range("A9").Value = '''''We love Taiwan'': Domestic workers hope for more from new ...'
$ws.Range("B9").Value = 45432
$ws.Range("C9").Value = 'While Taiwan''s monthly minimum salary was increased to 27,470 New Taiwan dollars ($853) this year, migrant domestic workers, who also have to'
$ws.Range("D9").Value = './Output/Images/7.jpg'
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = $true

$ws.Range("A10").Value = 'Who died alongside Iran''s President Raisi in the helicopter crash ...'
$ws.Range("B10").Value = 45431
$ws.Range("C10").Value = 'The AQR is a colossal bonyad, or charitable trust, that has billions of dollars in assets and is the custodian of the shrine of Imam Reza'
$ws.Range("D10").Value = './Output/Images/8.jpg'
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = $false

$ws.Range("A11").Value = 'Iran helicopter crash updates: President Raisi, FM Amirabdollahian ...'
$ws.Range("B11").Value = 45431
$ws.Range("C11").Value = 'dollars, according to a Reuters investigation. Under Mokhber''s watch, Setad developed Iran''s coronavirus vaccine, Coviran Barekat, at the'
$ws.Range("D11").Value = './Output/Images/9.jpg'
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = $false

$ws.Range("A12").Value = 'Who is Mohammad Mokhber, Iran''s interim president? | Politics ...'
$ws.Range("B12").Value = 45431
$ws.Range("C12").Value = 'Mokhber led the Iranian supreme leader''s multibillion-dollar charitable conglomerate for 14 years.'
$ws.Range("D12").Value = './Output/Images/10.jpg'
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = $false

$ws.Range("A13").Value = 'Ebrahim Raisi, Iran''s president, dies in helicopter crash aged 63 ...'
$ws.Range("B13").Value = 45431
$ws.Range("C13").Value = 'The colossal bonyad, or charitable trust, has billions of dollars in assets and is the custodian of the shrine of Imam Reza, the eighth Shia'
$ws.Range("D13").Value = './Output/Images/11.jpg'
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = $false

$ws.Range("A14").Value = 'Panic in Bishkek: Why were Pakistani students attacked in ...'
$ws.Range("B14").Value = 45431
$ws.Range("C14").Value = 'The Pakistani rupee, which stood at 160 against the US dollar in December 2020, has since slipped by more than 70 percent to 278 rupees a dollar'
$ws.Range("D14").Value = './Output/Images/12.jpg'
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = $false

$ws.Range("A15").Value = 'Russian court seizes two European banks'' assets amid Western ...'
$ws.Range("B15").Value = 45430
$ws.Range("C15").Value = 'Freezing hundreds of billions of dollars in lenders'' assets was part of dispute over gas project halted by sanctions.'
$ws.Range("D15").Value = './Output/Images/13.jpg'
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = $false

$ws.Range("A16").Value = 'Lebanon''s economic crisis endures, as does the EU''s ''fear'' of ...'
$ws.Range("B16").Value = 45429
$ws.Range("C16").Value = 'Billions of dollars go to the Syrian government, leaving it at the centre of the amphetamine''s trade. Published On 10 Mar 202410 Mar 2024.'
$ws.Range("D16").Value = './Output/Images/14.jpg'
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = $false

$ws.Range("A17").Value = 'What is Trident, the US floating pier off Gaza? Will it work? | Israel ...'
$ws.Range("B17").Value = 45428
$ws.Range("C17").Value = 'Washington has provided billions of dollars in aid as well as weapons that Israel has used in Gaza since October 7. Source: Al Jazeera'
$ws.Range("D17").Value = './Output/Images/15.jpg'
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = $false

$ws.Range("A18").Value = 'Western volunteers join the battle against Myanmar''s military regime ...'
$ws.Range("B18").Value = 45428
$ws.Range("C18").Value = 'dollar arsenal supplied by Russia and China. Ethnic armies, public donations and weapon seizures partly as a result of last year''s Operation'
$ws.Range("D18").Value = './Output/Images/16.jpg'
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = $false

$ws.Range("A19").Value = 'Republicans in US House pass bill pushing Biden to send weapons ...'
$ws.Range("B19").Value = 45428
$ws.Range("C19").Value = 'Israel, a major recipient of US military assistance for decades, is still due to get billions of dollars of US weaponry, despite the delay'
$ws.Range("D19").Value = './Output/Images/17.jpg'
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = $false

$ws.Range("A20").Value = 'US announces $2bn in new aid for Ukraine as Russian forces ...'
$ws.Range("B20").Value = 45426
$ws.Range("C20").Value = 'Antony Blinken says US rushing military support as Ukraine struggles to hold off renewed Russian offensive.'
$ws.Range("D20").Value = './Output/Images/18.jpg'
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = $true

$ws.Range("A21").Value = 'Anxious Zimbabwean migrants, smugglers watch South Africa''s ...'
$ws.Range("B21").Value = 45426
$ws.Range("C21").Value = 'The border province even favours using the South African rand, which people prefer to the local currency or the US dollar, which is popular'
$ws.Range("D21").Value = './Output/Images/19.jpg'
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = $false

$ws.Range("A22").Value = 'Biden administration plans to send $1bn in military aid to Israel ...'
$ws.Range("B22").Value = 45426
$ws.Range("C22").Value = 'Request for tank ammunition, tactical vehicles for Israel despite Biden''s earlier pause on bombs over Rafah assault.'
$ws.Range("D22").Value = './Output/Images/20.jpg'
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = $true

$ws.Range("A23").Value = 'Russia''s defence rejig: ''Unfortunately for Ukraine, a very effective ...'
$ws.Range("B23").Value = 45425
$ws.Range("C23").Value = 'dollars on new weaponry and payments to servicemen and their families. “Putin needs an ''arsenal of autocracy'' that can outperform Ukraine'
$ws.Range("D23").Value = './Output/Images/21.jpg'
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = $false

$ws.Range("A24").Value = 'Lawrence Wong set to take centre stage as Singapore''s new prime ...'
$ws.Range("B24").Value = 45425
$ws.Range("C24").Value = 'dollars ($1.6m) a year including bonuses. “Wong''s biggest challenge in the short term will be to articulate an easy-to-understand, inclusive'
$ws.Range("D24").Value = './Output/Images/22.jpg'
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = $true

$ws.Range("A25").Value = 'Not even the US government knows the US government line on ...'
$ws.Range("B25").Value = 45425
$ws.Range("C25").Value = 'Of course, this money was authorised on top of the billions of dollars that the US already sends the country on an annual basis. When on May'
$ws.Range("D25").Value = './Output/Images/23.jpg'
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = $false

$ws.Range("A26").Value = 'Boeing''s jets turn 70: A timeline of highs, lows and turbulence ...'
$ws.Range("B26").Value = 45425
$ws.Range("C26").Value = 'That same model rocket would be used for the Apollo 11 mission in 1969, landing astronauts on the moon. Boeing, the billion-dollar company.'
$ws.Range("D26").Value = './Output/Images/24.jpg'
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = $false

$ws.Range("A27").Value = 'One of the biggest hurdles for athletes on the Olympic path: Money ...'
$ws.Range("B27").Value = 45424
$ws.Range("C27").Value = 'She declined to share the dollar figure for those costs as well. Lozano told Al Jazeera that she''s using the funds from her GoFundMe'
$ws.Range("D27").Value = './Output/Images/25.jpg'
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = $false

$ws.Range("A28").Value = 'Energy summit seeks to curb cooking habits that kill millions every ...'
$ws.Range("B28").Value = 45424
$ws.Range("C28").Value = 'dollars to fund expanded access to clean cooking methods.   “Dollar for dollar, it''s hard to imagine a single intervention that'
$ws.Range("D28").Value = './Output/Images/26.jpg'
$ws.Range("E28").Value = 3
$ws.Range("F28").Value = $false

$ws.Range("A29").Value = 'US university ties to weapons contractors under scrutiny amid war in ...'
$ws.Range("B29").Value = 45424
$ws.Range("C29").Value = 'Many student demonstrators have zeroed in on their schools'' multimillion-dollar endowment funds as a target for their activism. Those'
$ws.Range("D29").Value = './Output/Images/27.jpg'
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = $false

$ws.Range("A30").Value = 'How US Big Tech supports Israel''s AI-powered genocide and ...'
$ws.Range("B30").Value = 45423
$ws.Range("C30").Value = 'Nvidia, the trillion-dollar chip behemoth powering the AI revolution, has also announced it is expanding its already large R&D operations in'
$ws.Range("D30").Value = './Output/Images/28.jpg'
$ws.Range("E30").Value = 1
$ws.Range("F30").Value = $false

$ws.Range("A31").Value = 'What did Biden say about US arms transfers to Israel and what does ...'
$ws.Range("B31").Value = 45420
$ws.Range("C31").Value = 'The US sends Israel $3.8bn in military aid annually, and Congress recently approved billions of dollars in additional support for the country.'
$ws.Range("D31").Value = './Output/Images/29.jpg'
$ws.Range("E31").Value = 1
$ws.Range("F31").Value = $true

$ws.Range("A32").Value = 'After decades of decline, Air India is betting billions on a comeback ...'
$ws.Range("B32").Value = 45420
$ws.Range("C32").Value = 'In recent decades, India''s national airline came to be seen as a cautionary tale of decline as it racked up billions of dollars in losses and'
$ws.Range("D32").Value = './Output/Images/30.jpg'
$ws.Range("E32").Value = 1
$ws.Range("F32").Value = $false

$ws.Range("A33").Value = 'Should India take from the rich, give the poor? A new election ...'
$ws.Range("B33").Value = 45419
$ws.Range("C33").Value = 'Inequality was worsened over the past decade of Modi''s rule. India has 271 dollar billionaires, third behind only China and the US — and world''s'
$ws.Range("D33").Value = './Output/Images/31.jpg'
$ws.Range("E33").Value = 1
$ws.Range("F33").Value = $false

$ws.Range("A34").Value = 'Pentagon chief confirms US pause on weapons shipment to Israel ...'
$ws.Range("B34").Value = 45419
$ws.Range("C34").Value = 'Over the years, the United States has provided tens of billions of dollars in military aid to Israel.” ''Iron-clad'' support. The Biden'
$ws.Range("D34").Value = './Output/Images/32.jpg'
$ws.Range("E34").Value = 1
$ws.Range("F34").Value = $false

$ws.Range("A35").Value = 'Zimbabwe''s illegal forex dealers use WhatsApp to find clients, evade ...'
$ws.Range("B35").Value = 45418
$ws.Range("C35").Value = 'dollar is the preferred medium of exchange. Everyone from state utilities to street vendors accepts payment in US dollars. Because of the'
$ws.Range("D35").Value = './Output/Images/33.jpg'
$ws.Range("E35").Value = 2
$ws.Range("F35").Value = $false

$ws.Range("A36").Value = 'Boeing postpones launch of Starliner space capsule after technical ...'
$ws.Range("B36").Value = 45418
$ws.Range("C36").Value = 'NASA in 2014 awarded multibillion-dollar contracts to Boeing and SpaceX to develop space capsules for the space agency to ferry astronauts and'
$ws.Range("D36").Value = './Output/Images/34.jpg'
$ws.Range("E36").Value = 1
$ws.Range("F36").Value = $false

$ws.Range("A37").Value = 'Australia''s Qantas to pay $79m over ''ghost flights'' furore | Aviation ...'
$ws.Range("B37").Value = 45417
$ws.Range("C37").Value = 'Australia''s flagship airline Qantas has agreed to pay $120 million Australian dollars ($79m) to settle a lawsuit over the sale of tickets'
$ws.Range("D37").Value = './Output/Images/35.jpg'
$ws.Range("E37").Value = 1
$ws.Range("F37").Value = $true

# The sheet now only needs 37 data+header rows; remove the old trailing row 38.
$ws.Rows.Item(38).Delete()
